$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2278323.2
$ws.Range("J17").Value = 2278323.2
$ws.Range("L17").Value = 6834969.600000001
$ws.Range("N17").Value = -6835305.600000001

$ws.Range("H33").Value = 6801.6665
$ws.Range("I33").Value = 10154.3
$ws.Range("J33").Value = 96.40000000000001
$ws.Range("K33").Value = 10154.3
$ws.Range("L33").Value = 96.40000000000001
$ws.Range("M33").Value = -9925.299999999999
$ws.Range("N33").Value = -554.4

$ws.Range("H44").Value = 41853.848
$ws.Range("J44").Value = 41853.848
$ws.Range("L44").Value = 41853.848
$ws.Range("N44").Value = -42777.848

$ws.Range("H48").Value = 1375
$ws.Range("I48").Value = 500
$ws.Range("J48").Value = 1500
$ws.Range("K48").Value = 1500
$ws.Range("L48").Value = 4500
$ws.Range("M48").Value = -1208
$ws.Range("N48").Value = -5084

$ws.Range("H56").Value = 1375
$ws.Range("I56").Value = 500
$ws.Range("J56").Value = 1500
$ws.Range("K56").Value = 1500
$ws.Range("L56").Value = 4500
$ws.Range("M56").Value = -966
$ws.Range("N56").Value = -5568

$ws.Range("H64").Value = 3451.6155
$ws.Range("I64").Value = 3040
$ws.Range("J64").Value = 4499.364
$ws.Range("K64").Value = 3040
$ws.Range("L64").Value = 4499.364
$ws.Range("M64").Value = -2792
$ws.Range("N64").Value = -4995.364

$ws.Range("H67").Value = 3451.6155
$ws.Range("I67").Value = 3040
$ws.Range("J67").Value = 4499.364
$ws.Range("K67").Value = 3040
$ws.Range("L67").Value = 4499.364
$ws.Range("M67").Value = -2182
$ws.Range("N67").Value = -6215.364

$ws.Range("H69").Value = 3652.6191
$ws.Range("I69").Value = 3013
$ws.Range("J69").Value = 4234.091
$ws.Range("K69").Value = 9039
$ws.Range("L69").Value = 12702.273
$ws.Range("M69").Value = -8165
$ws.Range("N69").Value = -14450.273

$ws.Range("H72").Value = 3652.6191
$ws.Range("I72").Value = 3013
$ws.Range("J72").Value = 4234.091
$ws.Range("K72").Value = 27117
$ws.Range("L72").Value = 38106.819
$ws.Range("M72").Value = -22749
$ws.Range("N72").Value = -46842.819

$ws.Range("H76").Value = 3579.18
$ws.Range("I76").Value = 2974.3242
$ws.Range("J76").Value = 5300.6924
$ws.Range("K76").Value = 2974.3242
$ws.Range("L76").Value = 5300.6924
$ws.Range("M76").Value = -2659.3242
$ws.Range("N76").Value = -5930.6924

$ws.Range("H79").Value = 3579.18
$ws.Range("I79").Value = 2974.3242
$ws.Range("J79").Value = 5300.6924
$ws.Range("K79").Value = 2974.3242
$ws.Range("L79").Value = 5300.6924
$ws.Range("M79").Value = -1882.3242
$ws.Range("N79").Value = -7484.6924

$ws.Range("H112").Value = 1110.2745
$ws.Range("J112").Value = 1106.9584
$ws.Range("L112").Value = 3320.8752
$ws.Range("N112").Value = -5536.8752

$ws.Range("H138").Value = 2455.1685
$ws.Range("I138").Value = 1466.475
$ws.Range("J138").Value = 3174.2183
$ws.Range("K138").Value = 4399.424999999999
$ws.Range("L138").Value = 9522.6549
$ws.Range("M138").Value = 740.5750000000007
$ws.Range("N138").Value = -19802.6549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 27500
$ws.Range("I60").Value = 5000
$ws.Range("K60").Value = 5000
$ws.Range("M60").Value = -4267

$ws.Range("H61").Value = 2200
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -3424

$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2314
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -11568
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 2627.95
$ws.Range("I132").Value = 1701.5385
$ws.Range("J132").Value = 4348.4287
$ws.Range("K132").Value = 5104.6155
$ws.Range("L132").Value = 13045.2861
$ws.Range("M132").Value = -2574.6155
$ws.Range("N132").Value = -18105.2861

$ws.Range("H136").Value = 2200
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3027.0203
$ws.Range("I31").Value = 1112.5883
$ws.Range("J31").Value = 7226.4194
$ws.Range("K31").Value = 1112.5883
$ws.Range("L31").Value = 7226.4194
$ws.Range("M31").Value = -817.5882999999999
$ws.Range("N31").Value = -7816.4194

$ws.Range("H34").Value = 3027.0203
$ws.Range("I34").Value = 1112.5883
$ws.Range("J34").Value = 7226.4194
$ws.Range("K34").Value = 1112.5883
$ws.Range("L34").Value = 7226.4194
$ws.Range("M34").Value = -910.5882999999999
$ws.Range("N34").Value = -7630.4194

$ws.Range("H62").Value = 3656.4285
$ws.Range("I62").Value = 3765.8333
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 3765.8333
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -3141.8333
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 3656.4285
$ws.Range("I65").Value = 3765.8333
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 3765.8333
$ws.Range("L65").Value = 3000
$ws.Range("M65").Value = -15709.1665
$ws.Range("N65").Value = -21240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 702516.6
$ws.Range("I5").Value = 332.8
$ws.Range("K5").Value = 998.4000000000001
$ws.Range("M5").Value = -886.4000000000001

$ws.Range("H58").Value = 1500
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 6000
$ws.Range("N58").Value = -6256

$ws.Range("H75").Value = 821.3
$ws.Range("I75").Value = 653.25
$ws.Range("J75").Value = 933.3333
$ws.Range("K75").Value = 1959.75
$ws.Range("L75").Value = 2799.9999
$ws.Range("M75").Value = -961.75
$ws.Range("N75").Value = -4795.9999

$ws.Range("H78").Value = 821.3
$ws.Range("I78").Value = 653.25
$ws.Range("J78").Value = 933.3333
$ws.Range("K78").Value = 5879.25
$ws.Range("L78").Value = 8399.9997
$ws.Range("M78").Value = -887.25
$ws.Range("N78").Value = -18383.9997

$ws.Range("H113").Value = 474.875
$ws.Range("I113").Value = 457
$ws.Range("J113").Value = 492.75
$ws.Range("K113").Value = 1371
$ws.Range("L113").Value = 1478.25
$ws.Range("M113").Value = 799
$ws.Range("N113").Value = -5818.25

$ws.Range("H117").Value = 801.2
$ws.Range("I117").Value = 384.66666
$ws.Range("J117").Value = 1078.8889
$ws.Range("K117").Value = 1153.99998
$ws.Range("L117").Value = 3236.6667
$ws.Range("M117").Value = 2288.00002
$ws.Range("N117").Value = -10120.6667

$ws.Range("H122").Value = 924.8
$ws.Range("I122").Value = 504.35
$ws.Range("J122").Value = 2606.6
$ws.Range("K122").Value = 4539.150000000001
$ws.Range("L122").Value = 23459.4
$ws.Range("M122").Value = -2089.150000000001
$ws.Range("N122").Value = -28359.4

$ws.Range("H135").Value = 702516.6
$ws.Range("I135").Value = 332.8
$ws.Range("K135").Value = 2995.2
$ws.Range("M135").Value = -460.2000000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1952.4
$ws.Range("I80").Value = 1947.8948
$ws.Range("J80").Value = 1966.6666
$ws.Range("K80").Value = 1947.8948
$ws.Range("L80").Value = 1966.6666
$ws.Range("M80").Value = -949.8948
$ws.Range("N80").Value = -3962.6666

$ws.Range("H83").Value = 1952.4
$ws.Range("I83").Value = 1947.8948
$ws.Range("J83").Value = 1966.6666
$ws.Range("K83").Value = 9739.474
$ws.Range("L83").Value = 9833.333000000001
$ws.Range("M83").Value = -4747.474
$ws.Range("N83").Value = -19817.333

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H122").Value = 2311.5813
$ws.Range("I122").Value = 1817.9259
$ws.Range("J122").Value = 3144.625
$ws.Range("K122").Value = 5453.7777
$ws.Range("L122").Value = 9433.875
$ws.Range("M122").Value = -3003.7777
$ws.Range("N122").Value = -14333.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2061.0557
$ws.Range("I68").Value = 1681.4286
$ws.Range("J68").Value = 2302.6365
$ws.Range("K68").Value = 1681.4286
$ws.Range("L68").Value = 2302.6365
$ws.Range("M68").Value = -932.4286
$ws.Range("N68").Value = -3800.6365

$ws.Range("H71").Value = 2061.0557
$ws.Range("I71").Value = 1681.4286
$ws.Range("J71").Value = 2302.6365
$ws.Range("K71").Value = 8407.143
$ws.Range("L71").Value = 11513.1825
$ws.Range("M71").Value = -4663.143
$ws.Range("N71").Value = -19001.1825

$ws.Range("H132").Value = 12059.214
$ws.Range("I132").Value = 13392.368
$ws.Range("K132").Value = 40177.104
$ws.Range("M132").Value = -37647.104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 259526.05
$ws.Range("I136").Value = 371115.97
$ws.Range("J136").Value = 8448.75
$ws.Range("K136").Value = 1113347.91
$ws.Range("L136").Value = 25346.25
$ws.Range("M136").Value = -1110797.91
$ws.Range("N136").Value = -30446.25
